$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 89.666664
$ws.Range("J55").Value = 85
$ws.Range("L55").Value = 85
$ws.Range("N55").Value = -513

$ws.Range("H64").Value = 10197.6
$ws.Range("I64").Value = 7999.5
$ws.Range("K64").Value = 7999.5
$ws.Range("M64").Value = -7751.5

$ws.Range("H67").Value = 10197.6
$ws.Range("I67").Value = 7999.5
$ws.Range("K67").Value = 7999.5
$ws.Range("M67").Value = -7141.5

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = ""
$ws.Range("N108").Value = 0

$ws.Range("H132").Value = 4647.2
$ws.Range("I132").Value = 4190
$ws.Range("K132").Value = 12570
$ws.Range("M132").Value = -10040

$ws.Range("H137").Value = 2772.9
$ws.Range("I137").Value = 2772.9
$ws.Range("K137").Value = 8318.700000000001
$ws.Range("M137").Value = -5768.700000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3530.5374
$ws.Range("J45").Value = 4064.848
$ws.Range("L45").Value = 4064.848
$ws.Range("N45").Value = -4818.848

$ws.Range("H102").Value = 14259.083
$ws.Range("I102").Value = 14259.083
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 14259.083
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = ""
$ws.Range("N102").Value = -12637.083

$ws.Range("H110").Value = 2028.9166
$ws.Range("I110").Value = 1705.3334
$ws.Range("J110").Value = 2999.6667
$ws.Range("K110").Value = 1705.3334
$ws.Range("L110").Value = 2999.6667
$ws.Range("M110").Value = 339.6666
$ws.Range("N110").Value = -7089.6667

$ws.Range("H132").Value = 2231.76
$ws.Range("I132").Value = 1896.7142
$ws.Range("J132").Value = 3990.75
$ws.Range("K132").Value = 5690.142599999999
$ws.Range("L132").Value = 11972.25
$ws.Range("M132").Value = -3160.142599999999
$ws.Range("N132").Value = -17032.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3854
$ws.Range("I134").Value = 3854
$ws.Range("K134").Value = 11562
$ws.Range("M134").Value = -9027

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H82").Value = 45000
$ws.Range("J82").Value = 45000
$ws.Range("L82").Value = 45000
$ws.Range("N82").Value = -45722

$ws.Range("H85").Value = 45000
$ws.Range("J85").Value = 45000
$ws.Range("L85").Value = 45000
$ws.Range("N85").Value = -47496

$ws.Range("H99").Value = 10125.2
$ws.Range("I99").Value = 8596.666999999999
$ws.Range("J99").Value = 10780.286
$ws.Range("K99").Value = 8596.666999999999
$ws.Range("L99").Value = 10780.286
$ws.Range("M99").Value = -7098.666999999999
$ws.Range("N99").Value = -13776.286

$ws.Range("H105").Value = 5296.909
$ws.Range("I105").Value = 3826.7
$ws.Range("K105").Value = 3826.7
$ws.Range("M105").Value = -2079.7

$ws.Range("H126").Value = 10125.2
$ws.Range("I126").Value = 8596.666999999999
$ws.Range("J126").Value = 10780.286
$ws.Range("K126").Value = 25790.001
$ws.Range("L126").Value = 32340.858
$ws.Range("M126").Value = -23320.001
$ws.Range("N126").Value = -37280.858

$ws.Range("H127").Value = 125000
$ws.Range("J127").Value = 125000
$ws.Range("L127").Value = 125000
$ws.Range("N127").Value = -134920

$ws.Range("H132").Value = 2249.4443
$ws.Range("I132").Value = 2203.0667
$ws.Range("K132").Value = 6609.2001
$ws.Range("M132").Value = -4079.2001

$ws.Range("H134").Value = 1696.75
$ws.Range("I134").Value = 1696.75
$ws.Range("K134").Value = 5090.25
$ws.Range("M134").Value = -2555.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4389.8
$ws.Range("I3").Value = 4389.8
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 13169.4
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = -13057.4

$ws.Range("H5").Value = 639.7692
$ws.Range("J5").Value = 1265.3334
$ws.Range("L5").Value = 3796.0002
$ws.Range("N5").Value = -4020.0002

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").Value = ""

$ws.Range("H51").Value = 4470.75
$ws.Range("I51").Value = 2629.3333
$ws.Range("K51").Value = 7887.999899999999
$ws.Range("M51").Value = -7427.999899999999

$ws.Range("H54").Value = 7999.5
$ws.Range("J54").Value = 7999.5
$ws.Range("L54").Value = 23998.5
$ws.Range("N54").Value = -25116.5

$ws.Range("H92").Value = 677.1539
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 677.1539
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = ""
$ws.Range("M92").Value = 2031.4617
$ws.Range("N92").Value = -4527.4617

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = ""
$ws.Range("N101").Value = 0

$ws.Range("H104").Value = 4794
$ws.Range("I104").Value = 4794
$ws.Range("K104").Value = 14382
$ws.Range("M104").Value = -11761

$ws.Range("H121").Value = 100665.7
$ws.Range("I121").Value = 347.16666
$ws.Range("K121").Value = 1041.49998
$ws.Range("M121").Value = 268.5000199999999

$ws.Range("H132").Value = 492.66666
$ws.Range("I132").Value = 494
$ws.Range("K132").Value = 4446
$ws.Range("M132").Value = -1916

$ws.Range("H135").Value = 639.7692
$ws.Range("J135").Value = 1265.3334
$ws.Range("L135").Value = 11388.0006
$ws.Range("N135").Value = -16458.0006

$ws.Range("H137").Value = 4299.476
$ws.Range("I137").Value = 3568.3635
$ws.Range("K137").Value = 10705.0905
$ws.Range("M137").Value = -5605.0905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 15148547
$ws.Range("I11").Value = 16844334
$ws.Range("K11").Value = 16844334
$ws.Range("M11").Value = -16844195

$ws.Range("H53").Value = 26249.5
$ws.Range("I53").Value = 18333
$ws.Range("K53").Value = 18333
$ws.Range("M53").Value = -17702

$ws.Range("H132").Value = 5438.9546
$ws.Range("J132").Value = 4391.8
$ws.Range("L132").Value = 13175.4
$ws.Range("N132").Value = -18235.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 12503
$ws.Range("I23").Value = 12503
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 12503
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = ""
$ws.Range("N23").Value = -12273

$ws.Range("H31").Value = 5010.5
$ws.Range("J31").Value = 5262
$ws.Range("L31").Value = 5262
$ws.Range("N31").Value = -5758

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = ""

$ws.Range("H46").Value = 4733
$ws.Range("J46").Value = 4600
$ws.Range("L46").Value = 4600
$ws.Range("N46").Value = -4976

$ws.Range("H55").Value = 152.41667
$ws.Range("I55").Value = 176.77777
$ws.Range("J55").Value = 79.333336
$ws.Range("K55").Value = 176.77777
$ws.Range("L55").Value = 79.333336
$ws.Range("M55").Value = -3.777770000000004
$ws.Range("N55").Value = -425.333336

$ws.Range("H61").Value = 1349.5
$ws.Range("I61").Value = 1349.5
$ws.Range("K61").Value = 1349.5
$ws.Range("M61").Value = -1147.5

$ws.Range("H113").Value = 1349.5
$ws.Range("I113").Value = 1349.5
$ws.Range("K113").Value = 1349.5
$ws.Range("M113").Value = 820.5

$ws.Range("H132").Value = 11609.637
$ws.Range("I132").Value = 13623
$ws.Range("J132").Value = 2549.5
$ws.Range("K132").Value = 40869
$ws.Range("L132").Value = 7648.5
$ws.Range("M132").Value = -38339
$ws.Range("N132").Value = -12708.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 3025000
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676

$ws.Range("H132").Value = 7811.7915
$ws.Range("I132").Value = 5146.1875
$ws.Range("J132").Value = 13143
$ws.Range("K132").Value = 15438.5625
$ws.Range("L132").Value = 39429
$ws.Range("M132").Value = -12908.5625
$ws.Range("N132").Value = -44489
